# "- Na abertura inicial dos ficheiros, ... encontrámos a brilhante ideia de
# passar os ficheiros como argumentos para o resto das funções." ->
# "... passar os catálogos como argumentos para o resto das funções."
#
# The word "ficheiros" must become "catálogos", and in doing so the single
# run that holds the whole sentence has to be split into three runs (text
# before / the replaced word / text after), matching the author's edit.
#
# A plain Find/Replace keeps everything inside one run, so instead we
# locate the exact Range of the word, temporarily bracket it with a
# Bookmark (which forces Word to split the surrounding run into three
# pieces with no extra formatting), replace the text while that split is
# in effect, and only then delete the now-empty-use bookmark.

$d = $word.ActiveDocument

$sentence = "ficheiros como argumentos"
$fullText = $d.Content.Text
$start = $fullText.IndexOf($sentence)
if ($start -lt 0) {
    throw "Could not locate target sentence in document"
}

# Range over just the word "ficheiros" (9 characters) that needs replacing.
$wordRange = $d.Range($start, $start + 9)

# Bracket it with a bookmark: this splits the enclosing run into
# [before][ficheiros][after] without touching formatting.
$d.Bookmarks.Add("tmp_catalogos_split", $wordRange)

# Replace the text while the split (bookmark) is still in place, so the
# replacement lands cleanly in its own, isolated run.
$mid = $d.Range($start, $start + 9)
$mid.Text = "catálogos"

# Clean up the bookmark now that the run split has been captured; this
# does not re-merge the three runs.
$d.Bookmarks("tmp_catalogos_split").Delete()
